# Update flow set up
# - Insert two new log entries (rows) at the top of the Stream_seepage log,
#   add a missing date to the entry that was previously the first dated row,
#   and refresh the selection/dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stream_seepage")

# Make sure this is the active sheet (it already is tabSelected in the file).
$ws.Activate()

# Insert two new blank rows right after the header row (row 1).
# The first Insert() pushes everything from row 2 down to row 3, the second
# Insert() pushes it down again to row 4, leaving rows 2 and 3 free for the
# two new log entries.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Copy the date-formatted style (style index 2 in the original file, a
# m/d/yyyy number format) from an existing dated cell down onto the three
# affected date cells (A2, A3 are brand new; A4 belongs to the entry that
# previously had no date set).
$ws.Cells.Item(5, 1).Copy()
$ws.Range("A2:A4").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New shared strings must be created in the same order as in the authored
# file (193 = "Rerun with 20m...", 194 = "With 4x upscaling...",
# 195 = "With the proper connectivity..."), so populate D3, then C3, then C2.

# --- Row 3: new "Vertical Refinement" entry dated 2023-03-17 ---
$ws.Cells.Item(3, 4).Value = "Rerun with 20m (10 layers unconfined)"
$ws.Cells.Item(3, 3).Value = "With 4x upscaling and the adjustments made to the 8x model, the run time was 1 hr 18 min, 7 mxiter, 0.13% CME, NSE=0.41`nTaking a closer look at the water budget it is clear that the increased connectivity and higher conductivity led to much higher rates of groundwater outflow through the GHB. `n100 parallel runs of 4x upscale took 17 hours - except had issue where only top 10 m were unconfined (5 layers) instead of 20 m(10 layers). updated version took 15 hours`nThe same large increase in GHB outflow occured in the 8x upscaling showing that it really is the dominant forcing in a connected environment"
$ws.Cells.Item(3, 1).Value = 45002
$ws.Cells.Item(3, 2).Value = "Vertical Refinement"
$ws.Rows.Item(3).RowHeight = 115.2

# --- Row 2: new "Vertical Refinement" entry dated 2023-03-23 ---
$ws.Cells.Item(2, 3).Value = "With the proper connectivity set up, well pumping goes from 100,000 to 10,000 m3/d, 25,000 to 10,000 m3/d. GHB acually only decreases a little bit, less than I would expect."
$ws.Cells.Item(2, 1).Value = 45008
$ws.Cells.Item(2, 2).Value = "Vertical Refinement"
$ws.Rows.Item(2).RowHeight = 28.8
# Row 2 has no Task (column D) entry, unlike row 3 - remove the leftover
# blank formatted cell that Insert() propagated from the column style.
$ws.Cells.Item(2, 4).Clear()

# --- Row 4: previously the first (undated) "Vertical Refinement" entry ---
# It kept its existing Subject/Notes text but now gets an explicit date.
$ws.Cells.Item(4, 1).Value = 45001

# Update the active selection to match the authored state.
[void]$ws.Range("A3").Select()

Write-Host "edit applied"
